$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of the data block (row 10),
# pushing the existing rows 10-28 down to 11-29.
$ws.Rows.Item(10).Insert()

$ws.Cells.Item(10, 1).Value = 5
$ws.Cells.Item(10, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(10, 3).Value = "Maule"
$ws.Cells.Item(10, 4).Value = 44469
$ws.Cells.Item(10, 5).Value = 7
$ws.Cells.Item(10, 6).Value = 100112026
$ws.Cells.Item(10, 7).Value = "Haba"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 500
$ws.Cells.Item(10, 11).Value = 9000
$ws.Cells.Item(10, 12).Value = 9000
$ws.Cells.Item(10, 13).Value = 9000
$ws.Cells.Item(10, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 16).Value = 360
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"
